# Generate Report for Handoff
# Update status text and timestamps to reflect a fresh handoff generation,
# and narrow the "datetime" columns that used to be overly wide.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps refreshed for the new handoff generation ---
$wsOverview.Range("G2").Value = "2016-08-27 12:56:11"
$wsDeDe.Range("H2").Value = "2016-08-27 12:56:11"
$wsZhCn.Range("H2").Value = "2016-08-27 12:56:07"

# --- Narrow the zh-cn / de-de datetime columns on Overview (E:F) and the
#     Status columns widths used elsewhere (col C on each language sheet).
#     Excel's ColumnWidth is quantized to whole pixels (1/6-character
#     steps at the default Calibri 11 metrics), so 16.333... is the input
#     that lands closest to the authored 17.2159881591797 width. ---
$wsOverview.Range("E1").ColumnWidth = 16.333333333333336
$wsOverview.Range("F1").ColumnWidth = 16.333333333333336
$wsZhCn.Range("C1").ColumnWidth = 16.333333333333336
$wsDeDe.Range("C1").ColumnWidth = 16.333333333333336
